$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking columns (G,H,I,J,K) must be forced to Text format so Excel
# stores them as text strings (matching the source data which stores all
# values, including numbers, as text) instead of auto-converting to numbers.
$ws.Range("G14:K25").NumberFormat = "@"

$ws.Range("A14").Value = " Dubai (DSC)"
$ws.Range("B14").Value = " October 25 2020"
$ws.Range("C14").Value = "Super Kings won by 8 wickets (with 8 balls remaining)"
$ws.Range("D14").Value = "Royal Challengers Bangalore"
$ws.Range("E14").Value = "Chennai Super Kings"
$ws.Range("F14").Value = "Aaron Finch "
$ws.Range("G14").Value = "15"
$ws.Range("H14").Value = "11"
$ws.Range("I14").Value = "3"
$ws.Range("J14").Value = "0"
$ws.Range("K14").Value = "136.36"

$ws.Range("A15").Value = " Dubai (DSC)"
$ws.Range("B15").Value = " September 24 2020"
$ws.Range("C15").Value = "Kings XI won by 97 runs"
$ws.Range("D15").Value = "Royal Challengers Bangalore"
$ws.Range("E15").Value = "Kings XI Punjab"
$ws.Range("F15").Value = "Aaron Finch "
$ws.Range("G15").Value = "20"
$ws.Range("H15").Value = "21"
$ws.Range("I15").Value = "3"
$ws.Range("J15").Value = "0"
$ws.Range("K15").Value = "95.23"

$ws.Range("A16").Value = " Abu Dhabi"
$ws.Range("B16").Value = " November 06 2020"
$ws.Range("C16").Value = "Sunrisers won by 6 wickets (with 2 balls remaining)"
$ws.Range("D16").Value = "Royal Challengers Bangalore"
$ws.Range("E16").Value = "Sunrisers Hyderabad"
$ws.Range("F16").Value = "Aaron Finch "
$ws.Range("G16").Value = "32"
$ws.Range("H16").Value = "30"
$ws.Range("I16").Value = "3"
$ws.Range("J16").Value = "1"
$ws.Range("K16").Value = "106.66"

$ws.Range("A17").Value = " Sharjah"
$ws.Range("B17").Value = " October 15 2020"
$ws.Range("C17").Value = "Kings XI won by 8 wickets"
$ws.Range("D17").Value = "Royal Challengers Bangalore"
$ws.Range("E17").Value = "Kings XI Punjab"
$ws.Range("F17").Value = "Aaron Finch "
$ws.Range("G17").Value = "20"
$ws.Range("H17").Value = "18"
$ws.Range("I17").Value = "2"
$ws.Range("J17").Value = "1"
$ws.Range("K17").Value = "111.11"

$ws.Range("A18").Value = " Dubai (DSC)"
$ws.Range("B18").Value = " October 05 2020"
$ws.Range("C18").Value = "Capitals won by 59 runs"
$ws.Range("D18").Value = "Royal Challengers Bangalore"
$ws.Range("E18").Value = "Delhi Capitals"
$ws.Range("F18").Value = "Aaron Finch "
$ws.Range("G18").Value = "13"
$ws.Range("H18").Value = "14"
$ws.Range("I18").Value = "1"
$ws.Range("J18").Value = "0"
$ws.Range("K18").Value = "92.85"

$ws.Range("A19").Value = " Dubai (DSC)"
$ws.Range("B19").Value = " September 28 2020"
$ws.Range("C19").Value = "Match tied (RCB won the one-over eliminator)"
$ws.Range("D19").Value = "Royal Challengers Bangalore"
$ws.Range("E19").Value = "Mumbai Indians"
$ws.Range("F19").Value = "Aaron Finch "
$ws.Range("G19").Value = "52"
$ws.Range("H19").Value = "35"
$ws.Range("I19").Value = "7"
$ws.Range("J19").Value = "1"
$ws.Range("K19").Value = "148.57"

$ws.Range("A20").Value = " Dubai (DSC)"
$ws.Range("B20").Value = " September 21 2020"
$ws.Range("C20").Value = "RCB won by 10 runs"
$ws.Range("D20").Value = "Royal Challengers Bangalore"
$ws.Range("E20").Value = "Sunrisers Hyderabad"
$ws.Range("F20").Value = "Aaron Finch "
$ws.Range("G20").Value = "29"
$ws.Range("H20").Value = "27"
$ws.Range("I20").Value = "1"
$ws.Range("J20").Value = "2"
$ws.Range("K20").Value = "107.40"

$ws.Range("A21").Value = " Dubai (DSC)"
$ws.Range("B21").Value = " October 17 2020"
$ws.Range("C21").Value = "RCB won by 7 wickets (with 2 balls remaining)"
$ws.Range("D21").Value = "Royal Challengers Bangalore"
$ws.Range("E21").Value = "Rajasthan Royals"
$ws.Range("F21").Value = "Aaron Finch "
$ws.Range("G21").Value = "14"
$ws.Range("H21").Value = "11"
$ws.Range("I21").Value = "0"
$ws.Range("J21").Value = "2"
$ws.Range("K21").Value = "127.27"

$ws.Range("A22").Value = " Abu Dhabi"
$ws.Range("B22").Value = " October 21 2020"
$ws.Range("C22").Value = "RCB won by 8 wickets (with 39 balls remaining)"
$ws.Range("D22").Value = "Royal Challengers Bangalore"
$ws.Range("E22").Value = "Kolkata Knight Riders"
$ws.Range("F22").Value = "Aaron Finch "
$ws.Range("G22").Value = "16"
$ws.Range("H22").Value = "21"
$ws.Range("I22").Value = "2"
$ws.Range("J22").Value = "0"
$ws.Range("K22").Value = "76.19"

$ws.Range("A23").Value = " Sharjah"
$ws.Range("B23").Value = " October 12 2020"
$ws.Range("C23").Value = "RCB won by 82 runs"
$ws.Range("D23").Value = "Royal Challengers Bangalore"
$ws.Range("E23").Value = "Kolkata Knight Riders"
$ws.Range("F23").Value = "Aaron Finch "
$ws.Range("G23").Value = "47"
$ws.Range("H23").Value = "37"
$ws.Range("I23").Value = "4"
$ws.Range("J23").Value = "1"
$ws.Range("K23").Value = "127.02"

$ws.Range("A24").Value = " Dubai (DSC)"
$ws.Range("B24").Value = " October 10 2020"
$ws.Range("C24").Value = "RCB won by 37 runs"
$ws.Range("D24").Value = "Royal Challengers Bangalore"
$ws.Range("E24").Value = "Chennai Super Kings"
$ws.Range("F24").Value = "Aaron Finch "
$ws.Range("G24").Value = "2"
$ws.Range("H24").Value = "9"
$ws.Range("I24").Value = "0"
$ws.Range("J24").Value = "0"
$ws.Range("K24").Value = "22.22"

$ws.Range("A25").Value = " Abu Dhabi"
$ws.Range("B25").Value = " October 03 2020"
$ws.Range("C25").Value = "RCB won by 8 wickets (with 5 balls remaining)"
$ws.Range("D25").Value = "Royal Challengers Bangalore"
$ws.Range("E25").Value = "Rajasthan Royals"
$ws.Range("F25").Value = "Aaron Finch "
$ws.Range("G25").Value = "8"
$ws.Range("H25").Value = "7"
$ws.Range("I25").Value = "2"
$ws.Range("J25").Value = "0"
$ws.Range("K25").Value = "114.28"

